$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new price record at row 176 (Macroferia Regional de
# Talca - Kiwi), pushing the existing rows 176-210 down to 177-211.
$ws.Rows.Item(176).Insert()

$ws.Cells.Item(176, 1).Value  = 5
$ws.Cells.Item(176, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(176, 3).Value  = "Maule"
$ws.Cells.Item(176, 4).Value  = 44511
$ws.Cells.Item(176, 5).Value  = 7
$ws.Cells.Item(176, 6).Value  = "Fruta"
$ws.Cells.Item(176, 7).Value  = 100101
$ws.Cells.Item(176, 8).Value  = "Berries"
$ws.Cells.Item(176, 9).Value  = 100101007
$ws.Cells.Item(176, 10).Value = "Kiwi"
$ws.Cells.Item(176, 11).Value = "Hayward"
$ws.Cells.Item(176, 12).Value = "Primera"
$ws.Cells.Item(176, 13).Value = 150
$ws.Cells.Item(176, 14).Value = 11000
$ws.Cells.Item(176, 15).Value = 11000
$ws.Cells.Item(176, 16).Value = 11000
$ws.Cells.Item(176, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(176, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(176, 19).Value = 1100
$ws.Cells.Item(176, 20).Value = 10
